$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 only had columns B..K before. L, M, N, O are brand-new cells that
# need the header style (bold/centered/bordered, same as D1) copied in explicitly
# -- xlPasteFormats (-4122) so we copy formatting only, not the cached value.
$ws.Range('D1').Copy()
$ws.Range('L1').PasteSpecial(-4122)
$ws.Range('D1').Copy()
$ws.Range('M1').PasteSpecial(-4122)
$ws.Range('D1').Copy()
$ws.Range('N1').PasteSpecial(-4122)
$ws.Range('D1').Copy()
$ws.Range('O1').PasteSpecial(-4122)

# --- Header row (row 1) ---
$ws.Range('B1').Value = 'Query'
$ws.Range('C1').Value = 'Recall'
$ws.Range('D1').Value = 'Cosine Precision'
$ws.Range('E1').Value = 'Cosine Relevant'
$ws.Range('F1').Value = 'Cosine F2'
$ws.Range('G1').Value = 'Cluster Precision'
$ws.Range('H1').Value = 'Cluster Relevant'
$ws.Range('I1').Value = 'Cluster F2'
$ws.Range('J1').Value = 'MVEE Precision'
$ws.Range('K1').Value = 'MVEE Relevant'
$ws.Range('L1').Value = 'MVEE F2'
$ws.Range('M1').Value = 'Hull Precision'
$ws.Range('N1').Value = 'Hull Relevant'
$ws.Range('O1').Value = 'Hull F2'

# --- Row 2 ---
$ws.Range('B2').Value = '("process adaptation" OR "processes adaptation" OR "customization of processes" OR "software processes customization" OR "software process customization" OR "customizing software processes" OR "process definition" OR "processes definition" OR "process composition" OR "compose processes" OR "processes composition" OR "process tailoring" OR "processes tailoring" OR "tailing of processes" OR "process development" OR "processes development" OR "process engineering" OR "processes engineering" OR "process design" OR "software process modelling" OR "software process modelling" OR "process implementation" OR "managing processes") AND ("family of software process" OR "family of software processes" OR "families of software process" OR "families of software processes" OR "software process line" OR "software process lines" OR "software processes line" OR "software processes lines" OR "process-line" OR "process-lines" OR "processes-line" OR "processes-lines" OR "software process family" OR "software processes family" OR "software process families" OR "software processes families" OR "process-family" OR "processes-family" OR "process-families" OR "processes-families" OR "software process variability" OR "software process variabilities" OR "software processes variability" OR "software processes variabilities" OR "variabilities in software processes" OR "process domain engineering" OR "processes domain engineering" OR "process feature" OR "process features" OR "processes feature" OR "processes features" OR "process asset reuse")'
$ws.Range('C2').Value = 0.37
$ws.Range('D2').Value = 0.32
$ws.Range('E2').Value = 32
$ws.Range('F2').Value = 0.36
$ws.Range('G2').Value = 0.4
$ws.Range('H2').Value = 40
$ws.Range('I2').Value = 0.38
$ws.Range('J2').Value = 0.28
$ws.Range('K2').Value = 28
$ws.Range('L2').Value = 0.35
$ws.Range('M2').Value = 0.28
$ws.Range('N2').Value = 28
$ws.Range('O2').Value = 0.35

# --- Row 3 ---
$ws.Range('B3').Value = 'Software Process Line'
$ws.Range('C3').Value = 0.6
$ws.Range('D3').Value = 0.02
$ws.Range('E3').Value = 310
$ws.Range('F3').Value = 0.1
$ws.Range('G3').Value = 0.08
$ws.Range('H3').Value = 1068
$ws.Range('I3').Value = 0.24
$ws.Range('J3').Value = 0.11
$ws.Range('K3').Value = 1523
$ws.Range('L3').Value = 0.31
$ws.Range('M3').Value = 0.08
$ws.Range('N3').Value = 1145
$ws.Range('O3').Value = 0.26

# --- Row 4 ---
$ws.Range('B4').Value = '(stream processing OR "continuous query" OR "stream-based system" OR "data stream system" OR "streaming system" OR "complex event processing") AND ("adapt" OR "reconfigur" ) AND ("latency" OR "response time")'
$ws.Range('C4').Value = 0.1
$ws.Range('D4').Value = 0.3
$ws.Range('E4').Value = 78
$ws.Range('F4').Value = 0.12
$ws.Range('G4').Value = 0.17
$ws.Range('H4').Value = 45
$ws.Range('I4').Value = 0.11
$ws.Range('J4').Value = 0.24
$ws.Range('K4').Value = 62
$ws.Range('L4').Value = 0.14
$ws.Range('M4').Value = 0.18
$ws.Range('N4').Value = 48
$ws.Range('O4').Value = 0.14

# --- Row 5 ---
$ws.Range('B5').Value = 'Data Stream Processing Latency'
$ws.Range('C5').Value = 0.17
$ws.Range('D5').Value = 0.3
$ws.Range('E5').Value = 528
$ws.Range('F5').Value = 0.19
$ws.Range('G5').Value = 0.49
$ws.Range('H5').Value = 865
$ws.Range('I5').Value = 0.21
$ws.Range('J5').Value = 0.39
$ws.Range('K5').Value = 683
$ws.Range('L5').Value = 0.21
$ws.Range('M5').Value = 0.31
$ws.Range('N5').Value = 550
$ws.Range('O5').Value = 0.2

# --- Row 6 ---
$ws.Range('B6').Value = '("metamodel" OR "meta-model") AND ("business process" OR "process model" OR "petrinet" OR "petri-net" OR "workflow" OR "Declare")'
$ws.Range('C6').Value = 0.58
$ws.Range('D6').Value = 0.3
$ws.Range('E6').Value = 380
$ws.Range('F6').Value = 0.49
$ws.Range('G6').Value = 0.26
$ws.Range('H6').Value = 326
$ws.Range('I6').Value = 0.4
$ws.Range('J6').Value = 0.56
$ws.Range('K6').Value = 705
$ws.Range('L6').Value = 0.63
$ws.Range('M6').Value = 0.43
$ws.Range('N6').Value = 549
$ws.Range('O6').Value = 0.59

# --- Row 7 ---
$ws.Range('B7').Value = 'Business Process Meta Models'
$ws.Range('C7').Value = 0.31
$ws.Range('D7').Value = 0.26
$ws.Range('E7').Value = 254
$ws.Range('F7').Value = 0.3
$ws.Range('G7').Value = 0.34
$ws.Range('H7').Value = 331
$ws.Range('I7').Value = 0.31
$ws.Range('J7').Value = 0.15
$ws.Range('K7').Value = 145
$ws.Range('L7').Value = 0.25
$ws.Range('M7').Value = 0.12
$ws.Range('N7').Value = 119
$ws.Range('O7').Value = 0.24

# --- Row 8 ---
$ws.Range('B8').Value = ' (("Parallel Programming") AND Modeling) OR (Multicore AND (Modeling OR "Software Engineering")) OR (Multicore AND ("Parallel Programming")) AND ("Modeling" OR "Software Engineering")'
$ws.Range('C8').Value = 0
$ws.Range('D8').Value = 0.57
$ws.Range('E8').Value = 31
$ws.Range('F8').Value = 0
$ws.Range('G8').Value = 0
$ws.Range('H8').Value = 0
$ws.Range('I8').Value = 0
$ws.Range('J8').Value = 0
$ws.Range('K8').Value = 0
$ws.Range('L8').Value = 0
$ws.Range('M8').Value = 0
$ws.Range('N8').Value = 0
$ws.Range('O8').Value = 0

# --- Row 9 ---
$ws.Range('B9').Value = 'Multicore Performance Prediction'
$ws.Range('C9').Value = 0
$ws.Range('D9').Value = 0.36
$ws.Range('E9').Value = 137
$ws.Range('F9').Value = 0
$ws.Range('G9').Value = 0
$ws.Range('H9').Value = 0
$ws.Range('I9').Value = 0
$ws.Range('J9').Value = 0
$ws.Range('K9').Value = 0
$ws.Range('L9').Value = 0
$ws.Range('M9').Value = 0
$ws.Range('N9').Value = 0
$ws.Range('O9').Value = 0

# --- Row 10 ---
$ws.Range('B10').Value = '((migration OR evolution OR adaptation OR transformation OR modernization OR reengineering OR integration OR adoption OR switching) AND (monolithic OR legacy OR existing OR preexisting OR on-premise) AND (system OR software OR application) AND (cloud AND (software OR application OR architecture OR infrastructure OR cloud environment)))'
$ws.Range('C10').Value = 0.24
$ws.Range('D10').Value = 0.08
$ws.Range('E10').Value = 177
$ws.Range('F10').Value = 0.17
$ws.Range('G10').Value = 0.04
$ws.Range('H10').Value = 90
$ws.Range('I10').Value = 0.12
$ws.Range('J10').Value = 0.01
$ws.Range('K10').Value = 28
$ws.Range('L10').Value = 0.05
$ws.Range('M10').Value = 0.01
$ws.Range('N10').Value = 17
$ws.Range('O10').Value = 0.03

# --- Row 11 ---
$ws.Range('B11').Value = 'Cloud Migration'
$ws.Range('C11').Value = 0.86
$ws.Range('D11').Value = 0.13
$ws.Range('E11').Value = 1098
$ws.Range('F11').Value = 0.4
$ws.Range('G11').Value = 0.04
$ws.Range('H11').Value = 296
$ws.Range('I11').Value = 0.15
$ws.Range('J11').Value = 0.15
$ws.Range('K11').Value = 1190
$ws.Range('L11').Value = 0.42
$ws.Range('M11').Value = 0.09
$ws.Range('N11').Value = 771
$ws.Range('O11').Value = 0.32

# --- Row 12 ---
$ws.Range('B12').Value = 'software AND (metric OR measurement) AND (fault OR defect* OR quality OR error-prone) AND (predict* OR prone* OR probability OR assess* OR detect* OR estimat* OR classificat*)'
$ws.Range('C12').Value = 0.81
$ws.Range('D12').Value = 0.09
$ws.Range('E12').Value = 1523
$ws.Range('F12').Value = 0.29
$ws.Range('G12').Value = 0.39
$ws.Range('H12').Value = 6906
$ws.Range('I12').Value = 0.54
$ws.Range('J12').Value = 0.17
$ws.Range('K12').Value = 2944
$ws.Range('L12').Value = 0.42
$ws.Range('M12').Value = 0.11
$ws.Range('N12').Value = 1870
$ws.Range('O12').Value = 0.33

# --- Row 13 ---
$ws.Range('B13').Value = 'Software Fault Prediction Metrics'
$ws.Range('C13').Value = 0.25
$ws.Range('D13').Value = 0.68
$ws.Range('E13').Value = 470
$ws.Range('F13').Value = 0.29
$ws.Range('G13').Value = 0.33
$ws.Range('H13').Value = 227
$ws.Range('I13').Value = 0.21
$ws.Range('J13').Value = 0.36
$ws.Range('K13').Value = 250
$ws.Range('L13').Value = 0.27
$ws.Range('M13').Value = 0.2
$ws.Range('N13').Value = 138
$ws.Range('O13').Value = 0.24

# --- Row 14 ---
$ws.Range('B14').Value = '(software OR applicati* OR systems ) AND (fault* OR defect* OR quality OR error-prone) AND (predict*OR prone* OR probability OR assess* OR detect* OR estimat* OR classificat*)'
$ws.Range('C14').Value = 0.2
$ws.Range('D14').Value = 0.02
$ws.Range('E14').Value = 355
$ws.Range('F14').Value = 0.08
$ws.Range('G14').Value = 0.02
$ws.Range('H14').Value = 295
$ws.Range('I14').Value = 0.07
$ws.Range('J14').Value = 0.02
$ws.Range('K14').Value = 332
$ws.Range('L14').Value = 0.08
$ws.Range('M14').Value = 0.01
$ws.Range('N14').Value = 153
$ws.Range('O14').Value = 0.04

# --- Row 15 ---
$ws.Range('B15').Value = 'Software Defect Prediction'
$ws.Range('C15').Value = 0.33
$ws.Range('D15').Value = 0.61
$ws.Range('E15').Value = 2170
$ws.Range('F15').Value = 0.36
$ws.Range('G15').Value = 0.47
$ws.Range('H15').Value = 1659
$ws.Range('I15').Value = 0.28
$ws.Range('J15').Value = 0.6
$ws.Range('K15').Value = 2124
$ws.Range('L15').Value = 0.36
$ws.Range('M15').Value = 0.46
$ws.Range('N15').Value = 1652
$ws.Range('O15').Value = 0.35
